$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CNXMOD -> context-cnxmod tags (B column) and add filter-type:multi-lo (D column)
$ws.Range("B2").Value = "context-cnxmod:39256206-03b0-4396-abb6-75e6ee5e3c7b"
$ws.Range("B3").Value = "context-cnxmod:102e9604-daa7-4a09-9f9e-232251d1a4ee"
$ws.Range("B4").Value = "alternate-context-cnxmod:39256206-03b0-4396-abb6-75e6ee5e3c7b,context-cnxmod:102e9604-daa7-4a09-9f9e-232251d1a4ee"
$ws.Range("D3").Value = "filter-type:multi-lo"
$ws.Range("D4").Value = "filter-type:multi-cnxmod,filter-type:multi-lo"

# Widen columns B and D to fit the new, longer tag text
$ws.Columns.Item(2).ColumnWidth = 90.42857142857143
$ws.Columns.Item(4).ColumnWidth = 29.142857142857142

# Move the active selection from A5 to D5
$ws.Range("D5").Select()
